# Apply the commit's changes to the "GD" sheet (Medida de generación distribuida):
# the "Inicio de Operaciones" (date) column was removed from the small summary
# table, so the "Electricidad generada (MWh)" and "Emisiones de GEI Reducidas"
# columns shift one column to the left. Deleting the whole column via COM
# automatically re-flows the remaining cells/styles/formulas/merged-cells and
# drops the now-unused "Inicio de Operaciones" shared string.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GD")

$ws.Columns("C").Delete() | Out-Null

# The GD sheet becomes the active / selected sheet & cell, replacing
# "Hoja3" (Proveedores' tab used to be tabSelected) as the active tab.
$ws.Activate()
$ws.Range("F11").Select() | Out-Null
